# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append this week's per-drive yardage logs to the running strings
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value2 + " 2 2 1 4 2 2 4 3 4 4 2 9 4 -3 1 10 6 2 0 4 4 1 5 3 4 0 0 5 -1 2"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value2 + " 3 4 3 1 -1 7 7 4 4 0 2 2 10 3 1 2 19 0 2 1 2 4"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value2 + " 4 4 5 -1 18 7 11 6 7 6 10 9 17 13 40 24 1 5 12"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value2 + " 4 17 8 11 12 4 14 -5 6 4 4 56"

# ---------------------------------------------------------------------------
# OFF sheet: updated counting stats after logging Week 16 + sim from Week 17
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 172
$wsOFF.Range("D2").Value = 9
$wsOFF.Range("F2").Value = 64
$wsOFF.Range("G2").Value = 48
$wsOFF.Range("H2").Value = 3
$wsOFF.Range("I2").Value = 6
$wsOFF.Range("J2").Value = 28
$wsOFF.Range("N2").Value = 12

$wsOFF.Range("B3").Value = 11
$wsOFF.Range("C3").Value = 171
$wsOFF.Range("E3").Value = 26
$wsOFF.Range("F3").Value = 103
$wsOFF.Range("G3").Value = 41
$wsOFF.Range("I3").Value = 54
$wsOFF.Range("J3").Value = 68
$wsOFF.Range("L3").Value = 259
$wsOFF.Range("M3").Value = 175
$wsOFF.Range("Q3").Value = 437

# ---------------------------------------------------------------------------
# DEF sheet: updated counting stats after logging Week 16 + sim from Week 17
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 172
$wsDEF.Range("F2").Value = 56
$wsDEF.Range("G2").Value = 28
$wsDEF.Range("H2").Value = 6
$wsDEF.Range("J2").Value = 23
$wsDEF.Range("N2").Value = 18
$wsDEF.Range("O2").Value = 17
$wsDEF.Range("P2").Value = 10

$wsDEF.Range("C3").Value = 180
$wsDEF.Range("E3").Value = 47
$wsDEF.Range("F3").Value = 113
$wsDEF.Range("H3").Value = 38
$wsDEF.Range("I3").Value = 66
$wsDEF.Range("J3").Value = 45
$wsDEF.Range("L3").Value = 267
$wsDEF.Range("M3").Value = 176
$wsDEF.Range("Q3").Value = 447

# ---------------------------------------------------------------------------
# ST sheet: special-teams counts + per-kick logs
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value = 70
$wsST.Range("D2").Value = 69
$wsST.Range("F2").Value = 134
$wsST.Range("G2").Value = 131
$wsST.Range("J2").Value = 55
$wsST.Range("K2").Value = 52
$wsST.Range("L2").Value = 41
$wsST.Range("M2").Value = 31
$wsST.Range("N2").Value = 20

$wsST.Range("B3").Value = 50

$wsST.Range("B4").Value = $wsST.Range("B4").Value2 + " 71"
$wsST.Range("B5").Value = $wsST.Range("B5").Value2 + " 26"
$wsST.Range("B6").Value = $wsST.Range("B6").Value2 + " 23"
$wsST.Range("D3").Value = $wsST.Range("D3").Value2 + " 41 37 47 43 42 45"
$wsST.Range("D4").Value = $wsST.Range("D4").Value2 + " 10 0 0 0 0 17"
$wsST.Range("D5").Value = $wsST.Range("D5").Value2 + " -1 9 0 0 0 0 10"

# ---------------------------------------------------------------------------
# TURNS sheet: turnover counts
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("B3").Value = 6
$wsTURNS.Range("C3").Value = 3
$wsTURNS.Range("D3").Value = 11
$wsTURNS.Range("E3").Value = 7

# ---------------------------------------------------------------------------
# PEN sheet: penalty counts
# ---------------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B3").Value = 17
